# Applies the cryptos.xlsx update described by the commit diff.
# Most cells are plain text (coin name / link / price-as-text / formatted pct)
# and are set via .Value. A handful of "Price" cells are strings that LOOK
# like plain numbers (e.g. "14.80"); assigning those through .Value lets Excel
# auto-coerce them to a Number and silently drop the trailing zero (14.80 -> 14.8).
# For those we force text entry with a leading apostrophe via .Formula, then
# reset .Style back to "Normal" so we do not leave a stray quote-prefixed style
# behind (matches the source workbook, which never touches cell styles).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '64.364.23'
$ws.Cells.Item(2, 5).Value = '  -3.30%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.143.72'
$ws.Cells.Item(3, 5).Value = '  -2.30%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '0.998'
$ws.Cells.Item(4, 5).Value = '  -0.08%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '602.15'
$ws.Cells.Item(5, 5).Value = '  -1.13%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '150.64'
$ws.Cells.Item(6, 5).Value = '  -4.83%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.05%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '3.141.90'
$ws.Cells.Item(8, 5).Value = '  -2.33%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.536'
$ws.Cells.Item(9, 5).Value = '  -2.84%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -4.19%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '5.62'
$ws.Cells.Item(11, 5).Value = '  -1.14%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -4.36%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -3.08%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '37.13'
$ws.Cells.Item(14, 5).Value = '  -4.14%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.617.38'
$ws.Cells.Item(15, 5).Value = '  -3.46%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '64.404.47'
$ws.Cells.Item(16, 5).Value = '  -3.33%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +0.42%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.137.83'
$ws.Cells.Item(18, 5).Value = '  -2.50%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -4.05%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '484.63'
$ws.Cells.Item(20, 5).Value = '  -4.25%  '

# Row 21 (Price "14.80" would be auto-coerced to the Number 14.8 and lose its
# trailing zero if assigned via .Value, so force literal text entry instead)
$c = $ws.Cells.Item(21, 4)
$c.Formula = "'14.80"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -2.37%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '0.717'
$ws.Cells.Item(22, 5).Value = '  -2.21%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -2.14%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '13.98'
$ws.Cells.Item(24, 5).Value = '  -4.30%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '84.99'
$ws.Cells.Item(25, 5).Value = '  +0.14%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -0.16%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '2.94'
$ws.Cells.Item(27, 5).Value = '  -2.07%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '8.72'
$ws.Cells.Item(28, 5).Value = '  -4.29%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Hedera'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(29, 4).Value = '0.127'
$ws.Cells.Item(29, 5).Value = '  +3.24%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'ImmutableX'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(30, 4).Value = '2.27'
$ws.Cells.Item(30, 5).Value = '  -4.00%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '7.15'
$ws.Cells.Item(31, 5).Value = '  +1.97%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '2.73'
$ws.Cells.Item(32, 5).Value = '  -7.10%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'EthereumClassic'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(33, 4).Value = '26.92'
$ws.Cells.Item(33, 5).Value = '  -4.23%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(34, 4).Value = '0.999'
$ws.Cells.Item(34, 5).Value = '  -0.22%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -6.21%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'dogwifhat'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(36, 4).Value = '3.34'
$ws.Cells.Item(36, 5).Value = '  +9.62%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'Filecoin'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(37, 4).Value = '6.14'
$ws.Cells.Item(37, 5).Value = '  -5.28%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '54.69'
$ws.Cells.Item(38, 5).Value = '  -1.29%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.0₃0759'
$ws.Cells.Item(39, 5).Value = '  -1.52%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '453.89'
$ws.Cells.Item(40, 5).Value = '  -9.69%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -4.16%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '0.0405'
$ws.Cells.Item(42, 5).Value = '  -3.92%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '8.57'
$ws.Cells.Item(43, 5).Value = '  -1.81%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '2.45'
$ws.Cells.Item(44, 5).Value = '  -0.27%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '2.904.95'
$ws.Cells.Item(45, 5).Value = '  -0.05%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '0.277'
$ws.Cells.Item(46, 5).Value = '  -6.59%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '27.14'
$ws.Cells.Item(47, 5).Value = '  -3.71%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +0.16%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '2.35'
$ws.Cells.Item(50, 5).Value = '  -2.48%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +0.36%  '
